# Append new trading-log rows (36-41) to the active worksheet, mirroring the
# rows written by the trading bot's logger on 2025-09-23.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36: TRADING_ATTEMPT for ARB
$ws.Range("A36").Value = "2025-09-23T01:25:25.131648"
$ws.Range("B36").Value = "TRADING_ATTEMPT"
$ws.Range("C36").Value = "ARB"
$ws.Range("D36").Value = "UNKNOWN"
$ws.Range("E36").Value = 0.4364983524810095
$ws.Range("K36").Value = "ATTEMPT"
$ws.Range("L36").Value = "Attempting trade 1/3"

# Row 37: POSITION_OPENED for ARB
$ws.Range("A37").Value = "2025-09-23T01:25:27.195967"
$ws.Range("B37").Value = "POSITION_OPENED"
$ws.Range("C37").Value = "ARB"
$ws.Range("D37").Value = "UNKNOWN"
$ws.Range("E37").Value = 0.4364983524810095
$ws.Range("F37").Value = 900
$ws.Range("G37").Value = 10
$ws.Range("H37").Value = 0.125864799809103
$ws.Range("K37").Value = "SUCCESS"

# Row 38: TRADING_ATTEMPT for XRP
$ws.Range("A38").Value = "2025-09-23T01:25:27.217209"
$ws.Range("B38").Value = "TRADING_ATTEMPT"
$ws.Range("C38").Value = "XRP"
$ws.Range("D38").Value = "UNKNOWN"
$ws.Range("E38").Value = 2.851139194678168
$ws.Range("K38").Value = "ATTEMPT"
$ws.Range("L38").Value = "Attempting trade 2/3"

# Row 39: POSITION_OPENED for XRP
$ws.Range("A39").Value = "2025-09-23T01:25:28.902338"
$ws.Range("B39").Value = "POSITION_OPENED"
$ws.Range("C39").Value = "XRP"
$ws.Range("D39").Value = "UNKNOWN"
$ws.Range("E39").Value = 2.851139194678168
$ws.Range("F39").Value = 1800
$ws.Range("G39").Value = 20
$ws.Range("H39").Value = 0.1159561015951102
$ws.Range("K39").Value = "SUCCESS"

# Row 40: TRADING_ATTEMPT for ENA
$ws.Range("A40").Value = "2025-09-23T01:25:28.923622"
$ws.Range("B40").Value = "TRADING_ATTEMPT"
$ws.Range("C40").Value = "ENA"
$ws.Range("D40").Value = "UNKNOWN"
$ws.Range("E40").Value = 0.6000566049103071
$ws.Range("K40").Value = "ATTEMPT"
$ws.Range("L40").Value = "Attempting trade 3/3"

# Row 41: POSITION_FAILED for ENA
$ws.Range("A41").Value = "2025-09-23T01:25:29.117100"
$ws.Range("B41").Value = "POSITION_FAILED"
$ws.Range("C41").Value = "ENA"
$ws.Range("D41").Value = "UNKNOWN"
$ws.Range("K41").Value = "FAILED"
$ws.Range("L41").Value = "Trade execution failed for trade 3"
